$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GPIO")

$xlPasteFormats = -4122

# --- Header row -----------------------------------------------------
# I1 wording changes from "Code" to "Definitions"; two new header cells
# are added for the "Global variables" (K1) and "Map" (M1) columns, using
# the same header formatting as the existing I1 header cell.
$ws.Range("I1").Value = "Definitions"

$ws.Range("I1").Copy()
$ws.Range("K1").PasteSpecial($xlPasteFormats)
$ws.Range("M1").PasteSpecial($xlPasteFormats)
$ws.Range("K1").Value = "Global variables"
$ws.Range("M1").Value = "Map"

# --- Data rows --------------------------------------------------------
# GPIO definition rows (row 8 and 15 are blank separator rows).
$rows = @(2,3,4,5,6,7,9,10,11,12,13,14,16,17,18,19,20,21,22,23)

foreach ($r in $rows) {
    $iCell = $ws.Cells.Item($r, 9)    # column I (existing "Definitions" formula cell, used as format source)
    $kCell = $ws.Cells.Item($r, 11)   # column K - "Global variables"
    $mCell = $ws.Cells.Item($r, 13)   # column M - "Map"

    $iCell.Copy()
    $kCell.PasteSpecial($xlPasteFormats)
    $mCell.PasteSpecial($xlPasteFormats)

    $kCell.Formula = "=CONCATENATE(`"static const px_gpio_handle_t `",LOWER(`$A$r),`" = {`",`$A$r,`"};`")"
    $mCell.Formula = "=CONCATENATE(`"#define PX_GPIO_`",`$B$r,`$C$r,`" `",`$A$r)"
}

# --- Column widths ------------------------------------------------------
# J and L stay default width; K and M get bestFit widths matching content.
$ws.Columns.Item(11).ColumnWidth = 54.7109375
$ws.Columns.Item(13).ColumnWidth = 33.85546875
